# Simulated Wild Card round and logged it
# Appends one game (Wild Card playoff) worth of stats to the cumulative
# per-play series (YDS, ST sheets) and to the summary totals
# (OFF, DEF, ST, TURNS, PEN sheets).

$wb = $excel.ActiveWorkbook

# ---- Long per-play / per-game numeric series (stored as space separated text) ----
$longSeries = @{
    "YDS!B2" = "2 15 4 -1 8 3 -1 -1 2 2 1 15 8 4 -3 5 2 5 4 2 6 18 -1 2 -1 6 2 1 3 2 8 3 4 -2 6 1 12 1 6 5 2 2 3 24 4 1 6 0 1 14 10 1 3 5 3 19 4 4 0 -1 0 1 15 4 2 4 2 0 6 3 -2 4 -1 4 3 2 0 0 4 2 4 0 3 1 13 7 4 5 -1 3 7 4 5 2 16 4 3 4 1 8 1 7 5 0 1 3 4 4 3 3 6 5 2 8 3 7 5 2 2 3 8 6 0 16 3 6 2 3 9 5 3 9 3 10 -5 4 1 7 5 5 5 -3 2 2 5 3 -2 12 3 1 1 0 2 6 2 0 3 1 7 0 1 2 4 4 6 -2 6 0 5 2 0 -2 7 2 2 2 5 -5 0 4 3 0 -1 2 3 1 -1 5 4 2 4 2 1 -1 4 1 14 0 11 0 1 7 2 1 14 3 -1 0 1 -3 5 0 4 1 2 7 3 2 20 9 9 0 -3 6 4 0 2 2 4 3 4 6 3 3 6 0 0 3 2 6 3 1 5 0 2 2 3 2 5 -3 3 0 2 16 0 2 4 4 3 3 8 7 2 0 8 2 1 9 2 1 1 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2 3 4 0 0 4 10 1 1 12 4 5 0 8 2 5 4 1 6 6 9 3 7 2 0 11 3 2 0 3 24 8 30 4 12 1 4 2 2 1 4 2 2 4 3 4 4 2 9 4 -3 1 10 6 2 0 4 4 1 5 3 4 0 0 5 -1 2 16 0 1 0 1 14 11 1 0 6 1 8 1 6 5 3 4 2 5 4 3 0 7 5 2 1 10 6 5 3 1 3 3 11 9 1 27 7 3 5 1 3 1 4 3 2 2 23 4 -1 5 2 -2 4 11 3 2 0"
    "YDS!B3" = "17 5 9 18 23 6 36 4 0 30 12 3 -2 15 3 13 13 7 6 9 5 18 12 4 5 14 4 11 16 11 3 5 8 0 4 2 7 8 4 5 6 8 5 -1 12 8 4 10 9 7 4 13 4 4 2 7 8 9 2 2 3 5 9 6 7 4 10 3 7 1 5 2 18 4 10 8 6 4 -7 11 8 7 4 22 25 42 22 4 21 4 34 9 12 10 5 24 5 3 12 5 2 21 6 10 8 4 23 10 2 7 1 0 11 5 21 7 10 17 4 20 6 6 6 22 25 15 3 2 8 1 5 20 7 8 5 14 0 6 32 20 -1 12 18 2 12 9 8 5 -1 9 17 5 2 10 10 0 7 2 15 9 7 19 5 8 11 5 14 6 6 9 18 27 9 11 1 7 10 17 16 4 23 7 8 0 14 15 22 12 8 3 8 5 1 -1 6 12 10 4 10 -2 8 17 18 12 3 6 12 10 5 11 18 21 5 6 10 9 1 4 9 11 2 13 14 6 6 7 16 9 11 7 21 2 6 19 2 52 11 9 1 6 14 35 23 6 64 3 26 1 19 7 5 8 13 3 4 7 7 6 9 1 8 65 15 8 0 8 1 7 6 5 7 6 2 5 -1 25 4 10 6 0 9 57 9 5 5 5 7 11 11 7 7 4 11 11 9 2 -2 5 6 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17 5 1 9 37 15 12 4 27 20 20 7 1 17 9 11 1 4 4 5 -1 18 7 11 6 7 6 10 9 17 13 40 24 1 5 12 20 0 8 14 14 11 15 21 0 10 6 6 45 1 2 9 8 9 4 5 13 7 9 7 17 4 3 5 9 2 0 15"
    "YDS!C2" = "35 3 6 1 2 2 11 4 5 2 0 3 4 4 1 4 5 1 -1 10 -4 1 3 3 1 8 1 3 5 2 46 4 6 1 10 1 -1 2 5 4 0 8 4 4 4 1 4 15 2 5 0 0 7 3 4 2 8 1 3 2 3 5 3 0 4 2 6 4 2 5 2 14 7 4 1 19 3 2 0 1 4 -2 1 6 0 1 27 8 1 4 4 8 0 2 10 1 1 23 38 1 6 9 -1 0 12 -5 6 2 4 7 5 0 -5 6 0 -3 -1 1 10 1 17 4 8 4 4 4 6 4 -1 13 13 3 6 1 4 5 5 7 0 5 3 1 2 4 24 1 4 9 3 12 1 2 2 5 -1 1 7 3 13 1 2 7 2 2 2 2 1 3 11 2 2 3 5 3 2 4 -1 6 34 1 0 -1 -2 6 5 7 13 2 6 2 12 3 -1 4 0 7 4 8 3 9 3 3 1 1 10 2 4 8 1 0 2 4 1 9 4 19 4 12 2 8 4 2 -2 0 0 5 2 11 0 -4 1 3 1 11 5 7 7 5 39 -1 5 4 15 1 1 4 2 5 0 2 5 3 4 1 4 4 4 2 1 15 5 2 2 3 1 1 -3 0 5 -1 5 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3 6 8 11 1 2 2 3 1 5 3 5 7 4 1 -5 11 0 4 6 6 20 -2 2 1 3 4 3 1 -1 7 7 4 4 0 2 2 10 3 1 2 19 0 2 1 2 4 0 4 5 -3 3 14 0 4 4 5 4 -1 6 2 5 21 2 7 5 3 2 4 3 1 0 12 0 9 35 -1 2 4 0 -3 0 0 -2 1 39 3 6 7 -1 8 4 15 -1 3 5 1 7 12 9 2 4 6 19 0 4 4 6 5 -1 2 4 1"
    "YDS!C3" = "4 7 4 9 8 5 22 6 7 12 8 25 9 7 21 26 9 6 17 16 4 1 6 10 3 8 3 7 11 1 35 5 6 1 7 7 22 13 6 7 8 6 9 2 41 3 24 7 23 6 10 25 8 18 8 5 23 12 21 1 4 15 10 12 23 7 23 16 3 16 32 34 1 7 6 12 1 5 4 28 6 3 2 3 17 8 -2 3 14 10 1 3 41 11 24 20 9 26 13 10 10 18 6 15 62 6 10 10 11 12 3 4 11 5 12 6 14 12 16 4 15 7 34 19 9 22 24 1 19 4 6 8 9 11 20 10 28 11 29 11 11 21 7 14 10 7 10 20 12 -3 23 11 7 2 13 12 6 13 15 5 39 9 6 4 4 26 49 4 11 12 6 11 23 28 1 11 8 1 20 3 15 7 8 6 15 4 10 6 3 3 8 8 17 11 14 19 13 1 11 -1 20 8 4 3 23 9 9 15 6 7 5 10 6 9 8 15 4 39 15 4 10 3 5 10 7 4 7 4 11 15 7 5 6 13 3 11 9 9 11 0 8 7 8 4 8 20 5 8 30 6 9 15 9 2 5 5 29 14 20 9 5 22 5 4 9 23 5 5 2 62 12 8 11 8 15 64 4 13 4 7 18 24 13 17 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11 16 2 14 19 4 11 6 23 1 9 24 10 4 17 8 11 12 4 14 -5 6 4 4 56 25 2 4 1 5 5 16 3 7 52 6 3 15 2 5 35 7 5 13 23 17 6 2 28 6 32 5 1 18 39 11 8"
    "ST!B4" = "41 70 58 60 62 61 66 67 57 66 66 62 63 57 66 69 40 66 67 71 64 63"
    "ST!B5" = "17 23 15 21 23 0 23 27 5 24 21 19 18 18 21 29 20 20 31 26 24 24"
    "ST!B6" = "31 15 18 20 13 18 7 22 25 16 15 17 23 16 18 20 18"
    "ST!D3" = "40 52 37 51 48 47 30 36 61 39 34 50 32 41 46 53 56 39 35 36 44 30 37 32 59 49 47 58 43 59 58 55 44 49 39 33 50 36 38 40 38 49 44 40 52 47 53 41 57 40 37 44 35 46 53 50 65 39 37 46 51 52 56 41 37 47 43 42 45 37 36 42 32 44 63 57 45"
    "ST!D4" = "7 0 -1 14 20 0 0 0 13 7 0 2 0 0 0 0 -1 0 0 0 0 0 0 0 0 0 17 0 16 0 9 0 0 0 2 0 0 0 0 0 0 0 9 0 0 0 9 0 19 0 0 0 0 0 0 15 11 0 -4 0 0 11 12 10 0 0 0 0 17 9 8 11 0 11 0 9 0"
    "ST!D5" = "0 18 4 0 0 8 0 7 12 -1 14 0 0 0 0 0 0 0 16 0 8 0 -4 8 15 0 10 18 0 6 0 0 0 0 0 4 0 0 0 0 0 7 7 0 0 0 0 7 0 0 0 0 0 0 0 0 5 -1 9 0 0 0 0 10 15 0 0 0 0 0"
}

foreach ($key in $longSeries.Keys) {
    $parts = $key.Split("!")
    $sheetName = $parts[0]
    $cellRef = $parts[1]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $longSeries[$key]
}

# ---- Cumulative summary totals ----
$cellValues = @{
    "OFF!C2" = 198
    "OFF!D2" = 11
    "OFF!E2" = 13
    "OFF!F2" = 73
    "OFF!G2" = 56
    "OFF!I2" = 10
    "OFF!J2" = 34
    "OFF!L2" = 315
    "OFF!M2" = 210
    "OFF!O2" = 22
    "OFF!Q2" = 602
    "OFF!B3" = 12
    "OFF!C3" = 196
    "OFF!F3" = 119
    "OFF!G3" = 44
    "OFF!H3" = 29
    "OFF!I3" = 63
    "OFF!J3" = 74
    "OFF!N3" = 24
    "DEF!C2" = 201
    "DEF!D2" = 9
    "DEF!F2" = 69
    "DEF!G2" = 37
    "DEF!J2" = 31
    "DEF!L2" = 314
    "DEF!M2" = 180
    "DEF!Q2" = 554
    "DEF!C3" = 196
    "DEF!D3" = 5
    "DEF!E3" = 51
    "DEF!F3" = 123
    "DEF!G3" = 29
    "DEF!H3" = 40
    "DEF!I3" = 71
    "DEF!J3" = 52
    "DEF!N3" = 29
    "ST!B2" = 78
    "ST!D2" = 77
    "ST!F2" = 137
    "ST!G2" = 134
    "ST!J2" = 57
    "ST!K2" = 54
    "ST!L2" = 42
    "ST!M2" = 32
    "ST!B3" = 56
    "TURNS!C2" = 11
    "TURNS!E2" = 8
    "TURNS!D3" = 10
    "TURNS!E3" = 7
    "PEN!B2" = 21
    "PEN!D4" = 7
}

foreach ($key in $cellValues.Keys) {
    $parts = $key.Split("!")
    $sheetName = $parts[0]
    $cellRef = $parts[1]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $cellValues[$key]
}

Write-Host "Wild Card round simulation logged."
